$wb = $excel.ActiveWorkbook

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1422.8928
$ws.Range("I2").Value = 1513
$ws.Range("J2").Value = 1260.7
$ws.Range("K2").Value = 1513
$ws.Range("L2").Value = 1260.7
$ws.Range("M2").Value = -1400
$ws.Range("N2").Value = -1486.7
$ws.Range("H32").Value = 6968.6724
$ws.Range("I32").Value = 5550.6226
$ws.Range("K32").Value = 5550.6226
$ws.Range("M32").Value = -5263.6226
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H96").Value = 9672
$ws.Range("J96").Value = 9672
$ws.Range("L96").Value = 9672
$ws.Range("N96").Value = -15164
$ws.Range("H116").Value = 1422.8928
$ws.Range("I116").Value = 1513
$ws.Range("J116").Value = 1260.7
$ws.Range("K116").Value = 1513
$ws.Range("L116").Value = 1260.7
$ws.Range("M116").Value = 781
$ws.Range("N116").Value = -5848.7

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1422.8928
$ws.Range("I3").Value = 1513
$ws.Range("J3").Value = 1260.7
$ws.Range("K3").Value = 1513
$ws.Range("L3").Value = 1260.7
$ws.Range("M3").Value = -1399
$ws.Range("N3").Value = -1488.7

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2459953
$ws.Range("I58").Value = 4547268
$ws.Range("J58").Value = 4288.706
$ws.Range("K58").Value = 4547268
$ws.Range("L58").Value = 4288.706
$ws.Range("M58").Value = -4547065
$ws.Range("N58").Value = -4694.706
$ws.Range("H86").Value = 2277
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2277
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2277
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -4523
$ws.Range("H89").Value = 2277
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2277
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 11385
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -22617
$ws.Range("H94").Value = 1712.5
$ws.Range("I94").Value = 2007.2
$ws.Range("J94").Value = 1502
$ws.Range("K94").Value = 2007.2
$ws.Range("L94").Value = 1502
$ws.Range("M94").Value = -1556.2
$ws.Range("N94").Value = -2404
$ws.Range("H136").Value = 2459953
$ws.Range("I136").Value = 4547268
$ws.Range("J136").Value = 4288.706
$ws.Range("K136").Value = 13641804
$ws.Range("L136").Value = 12866.118
$ws.Range("M136").Value = -13639254
$ws.Range("N136").Value = -17966.118

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11907660
$ws.Range("I5").Value = 333.2258
$ws.Range("J5").Value = 45464670
$ws.Range("K5").Value = 999.6774
$ws.Range("L5").Value = 136394010
$ws.Range("M5").Value = -887.6774
$ws.Range("N5").Value = -136394234
$ws.Range("H113").Value = 701.9794000000001
$ws.Range("I113").Value = 737.2083
$ws.Range("J113").Value = 600.52
$ws.Range("K113").Value = 2211.6249
$ws.Range("L113").Value = 1801.56
$ws.Range("M113").Value = -41.6248999999998
$ws.Range("N113").Value = -6141.559999999999
$ws.Range("H118").Value = 7568.778
$ws.Range("I118").Value = 5950
$ws.Range("J118").Value = 8031.2856
$ws.Range("K118").Value = 17850
$ws.Range("L118").Value = 24093.8568
$ws.Range("M118").Value = -16607
$ws.Range("N118").Value = -26579.8568
$ws.Range("H135").Value = 11907660
$ws.Range("I135").Value = 333.2258
$ws.Range("J135").Value = 45464670
$ws.Range("K135").Value = 2999.0322
$ws.Range("L135").Value = 409182030
$ws.Range("M135").Value = -464.0322000000001
$ws.Range("N135").Value = -409187100

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 7666.6665
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888
$ws.Range("H8").Value = 7666.6665
$ws.Range("I8").Value = 3000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2861
$ws.Range("H80").Value = 5812.8184
$ws.Range("I80").Value = 18433
$ws.Range("J80").Value = 3820.158
$ws.Range("K80").Value = 18433
$ws.Range("L80").Value = 3820.158
$ws.Range("M80").Value = -17435
$ws.Range("N80").Value = -5816.157999999999
$ws.Range("H83").Value = 5812.8184
$ws.Range("I83").Value = 18433
$ws.Range("J83").Value = 3820.158
$ws.Range("K83").Value = 92165
$ws.Range("L83").Value = 19100.79
$ws.Range("M83").Value = -87173
$ws.Range("N83").Value = -29084.79

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 9500
$ws.Range("I3").Value = 4000
$ws.Range("J3").Value = 15000
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = -3888
$ws.Range("N3").Value = -15224
$ws.Range("H14").Value = 7940
$ws.Range("I14").Value = 4850
$ws.Range("K14").Value = 4850
$ws.Range("M14").Value = -4678
$ws.Range("H15").Value = 9500
$ws.Range("I15").Value = 4000
$ws.Range("J15").Value = 15000
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = -3830
$ws.Range("N15").Value = -15340
$ws.Range("H40").Value = 4392.32
$ws.Range("I40").Value = 4156
$ws.Range("K40").Value = 4156
$ws.Range("M40").Value = -4020
$ws.Range("H61").Value = 20591.46
$ws.Range("I61").Value = 51752
$ws.Range("J61").Value = 6742.3335
$ws.Range("K61").Value = 51752
$ws.Range("L61").Value = 6742.3335
$ws.Range("M61").Value = -51550
$ws.Range("N61").Value = -7146.3335
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null
$ws.Range("H113").Value = 20591.46
$ws.Range("I113").Value = 51752
$ws.Range("J113").Value = 6742.3335
$ws.Range("K113").Value = 51752
$ws.Range("L113").Value = 6742.3335
$ws.Range("M113").Value = -49582
$ws.Range("N113").Value = -11082.3335
$ws.Range("H136").Value = 5156.755
$ws.Range("I136").Value = 4044.6155
$ws.Range("K136").Value = 12133.8465
$ws.Range("M136").Value = -9583.8465

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 20333
$ws.Range("J12").Value = 20333
$ws.Range("L12").Value = 20333
$ws.Range("N12").Value = -20617
$ws.Range("H28").Value = 30000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 30000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = $null
$ws.Range("N28").Value = -30696
$ws.Range("H30").Value = 9822
$ws.Range("J30").Value = 8827.5
$ws.Range("L30").Value = 8827.5
$ws.Range("N30").Value = -9041.5
$ws.Range("H136").Value = 5632
$ws.Range("I136").Value = 3072.9473
$ws.Range("J136").Value = 7657.9165
$ws.Range("K136").Value = 9218.841899999999
$ws.Range("L136").Value = 22973.7495
$ws.Range("M136").Value = -6668.841899999999
$ws.Range("N136").Value = -28073.7495
